# Updates cryptocurrency price (column D) and 1h volume change (column E)
# values on the "cryptos" worksheet to the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.762.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "'2.491.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'586.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "'176.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  +3.32%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").Value = "'0.339"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value = "'4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'2.947.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "'25.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'67.672.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'2.404.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Value = "'10.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "'351.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'4.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'70.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.01%  "
$ws.Range("D24").Value = "'4.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "'1.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.31%  "
$ws.Range("D26").Value = "'9.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").Value = "'2.619.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "'0.0₃0904"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'504.50"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").Value = "'163.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "'144.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "'0.0742"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("D50").Value = "'0.585"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +0.45%  "
